$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.027.71'
$ws.Range('E2').Value = '  -3.50%  '

$ws.Range('D3').Value = '1.602.80'
$ws.Range('E3').Value = '  -2.74%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.001'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.06%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '301.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3776'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.87%  '

$ws.Range('E8').Value = '  -5.43%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '49.59'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.25%  '

$ws.Range('E10').Value = '  -6.02%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.001'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.09%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08117'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.74%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.82'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.08%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.584'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.03%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.403'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.65%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001246'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.01%  '

$ws.Range('D17').Value = '1.599.09'
$ws.Range('E17').Value = '  -2.99%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06877'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.23%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.49%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.563'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.46%  '

$ws.Range('E22').Value = '  +0.17%  '

$ws.Range('E23').Value = '  -3.38%  '

$ws.Range('D24').Value = '23.021.67'
$ws.Range('E24').Value = '  -3.55%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.366'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.18%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.804'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.56%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.08'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.91%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.36%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.258'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.27%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.82%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.299'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.56%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.824'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.45%  '

$ws.Range('D33').Value = '1.781.79'
$ws.Range('E33').Value = '  -2.66%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9641'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.40%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07628'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.13%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.35%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.299'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.00%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02706'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.15%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2534'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.34%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.08837'
$ws.Range('D40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.365'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.96%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7052'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.51%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.51'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.69%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.40%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6614'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.33%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.314'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.03%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.03%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.990'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.45%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.32%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07912'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.226'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.59%  '
